# Auto-generated script: applies scheduled market-data refresh to Leve profit sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 648.7778
$ws.Range("I2").Value = 540
$ws.Range("K2").Value = 540
$ws.Range("M2").Value = -427
$ws.Range("H62").Value = 3920.6
$ws.Range("I62").Value = 3934
$ws.Range("J62").Value = 3800
$ws.Range("K62").Value = 3934
$ws.Range("L62").Value = 3800
$ws.Range("M62").Value = -3310
$ws.Range("N62").Value = -5048
$ws.Range("H65").Value = 3920.6
$ws.Range("I65").Value = 3934
$ws.Range("J65").Value = 3800
$ws.Range("K65").Value = 19670
$ws.Range("L65").Value = 19000
$ws.Range("M65").Value = -16550
$ws.Range("N65").Value = -25240
$ws.Range("H70").Value = 5636.4
$ws.Range("I70").Value = 1996
$ws.Range("K70").Value = 5988
$ws.Range("M70").Value = -5718
$ws.Range("H73").Value = 5636.4
$ws.Range("I73").Value = 1996
$ws.Range("K73").Value = 5988
$ws.Range("M73").Value = -5052
$ws.Range("H74").Value = 14421.75
$ws.Range("I74").Value = 14716.533
$ws.Range("K74").Value = 14716.533
$ws.Range("M74").Value = -13780.533
$ws.Range("H77").Value = 14421.75
$ws.Range("I77").Value = 14716.533
$ws.Range("K77").Value = 73582.66499999999
$ws.Range("M77").Value = -68902.66499999999
$ws.Range("H112").Value = 1469.7142
$ws.Range("J112").Value = 1469.7142
$ws.Range("L112").Value = 4409.142599999999
$ws.Range("N112").Value = -6625.142599999999
$ws.Range("H116").Value = 4686.647
$ws.Range("J116").Value = 4394.7144
$ws.Range("L116").Value = 4394.7144
$ws.Range("N116").Value = -11278.7144

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1257.8
$ws.Range("I2").Value = 1170.1818
$ws.Range("K2").Value = 1170.1818
$ws.Range("M2").Value = -1057.1818
$ws.Range("H12").Value = 1251.5
$ws.Range("I12").Value = 1251.5
$ws.Range("K12").Value = 1251.5
$ws.Range("M12").Value = -1078.5
$ws.Range("H32").Value = 1193696.4
$ws.Range("I32").Value = 546075.75
$ws.Range("K32").Value = 546075.75
$ws.Range("M32").Value = -545788.75
$ws.Range("H74").Value = 2191.3416
$ws.Range("I74").Value = 1999.3871
$ws.Range("K74").Value = 1999.3871
$ws.Range("M74").Value = -1125.3871
$ws.Range("H77").Value = 2191.3416
$ws.Range("I77").Value = 1999.3871
$ws.Range("K77").Value = 9996.9355
$ws.Range("M77").Value = -5628.9355
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()  # was -32523
$ws.Range("H116").Value = 1257.8
$ws.Range("I116").Value = 1170.1818
$ws.Range("K116").Value = 1170.1818
$ws.Range("M116").Value = 1123.8182
$ws.Range("H125").Value = 124463.75
$ws.Range("J125").Value = 124463.75
$ws.Range("L125").Value = 124463.75
$ws.Range("N125").Value = -134303.75

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1257.8
$ws.Range("I3").Value = 1170.1818
$ws.Range("K3").Value = 1170.1818
$ws.Range("M3").Value = -1056.1818
$ws.Range("H107").Value = 3206309
$ws.Range("I107").Value = 4274455
$ws.Range("K107").Value = 4274455
$ws.Range("M107").Value = -4272535

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 371.6
$ws.Range("I7").Value = 120
$ws.Range("J7").Value = 623.2
$ws.Range("K7").Value = 120
$ws.Range("L7").Value = 623.2
$ws.Range("M7").Value = -7
$ws.Range("N7").Value = -849.2
$ws.Range("H31").Value = 3576482.2
$ws.Range("I31").Value = 3291.75
$ws.Range("K31").Value = 3291.75
$ws.Range("M31").Value = -2996.75
$ws.Range("H34").Value = 3576482.2
$ws.Range("I34").Value = 3291.75
$ws.Range("K34").Value = 3291.75
$ws.Range("M34").Value = -3089.75
$ws.Range("H116").Value = 70075.336
$ws.Range("J116").Value = 70075.336
$ws.Range("L116").Value = 70075.336
$ws.Range("N116").Value = -79253.336
$ws.Range("H120").Value = 49997
$ws.Range("J120").Value = 49997
$ws.Range("L120").Value = 49997
$ws.Range("N120").Value = -57255

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 4598
$ws.Range("J114").Value = 4665
$ws.Range("L114").Value = 13995
$ws.Range("N114").Value = -20503

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 40000
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()  # was -20771
$ws.Range("H122").Value = 4811662
$ws.Range("I122").Value = 6996462.5
$ws.Range("J122").Value = 5101.4
$ws.Range("K122").Value = 20989387.5
$ws.Range("L122").Value = 15304.2
$ws.Range("M122").Value = -20986937.5
$ws.Range("N122").Value = -20204.2
$ws.Range("H134").Value = 17000
$ws.Range("J134").Value = 17000
$ws.Range("L134").Value = 51000
$ws.Range("N134").Value = -56070
$ws.Range("H136").Value = 115999.8
$ws.Range("J136").Value = 115999.8
$ws.Range("L136").Value = 347999.4
$ws.Range("N136").Value = -353099.4

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 5000
$ws.Range("I58").Value = 5000
$ws.Range("K58").Value = 5000
$ws.Range("M58").Value = -4740
$ws.Range("H122").Value = 9610.223
$ws.Range("I122").Value = 7035.467
$ws.Range("J122").Value = 12828.667
$ws.Range("K122").Value = 21106.401
$ws.Range("L122").Value = 38486.001
$ws.Range("M122").Value = -18656.401
$ws.Range("N122").Value = -43386.001
$ws.Range("H132").Value = 13214
$ws.Range("J132").Value = 13374.75
$ws.Range("L132").Value = 40124.25
$ws.Range("N132").Value = -45184.25
$ws.Range("H136").Value = 2563.625
$ws.Range("I136").Value = 2001.2858
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 6003.857400000001
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = -3453.857400000001
$ws.Range("N136").Value = -24600

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 15873.333
$ws.Range("I58").Value = 15873.333
$ws.Range("K58").Value = 15873.333
$ws.Range("M58").Value = -15565.333
$ws.Range("H81").Value = 3278.0435
$ws.Range("J81").Value = 5519.6
$ws.Range("L81").Value = 11039.2
$ws.Range("N81").Value = -13161.2
$ws.Range("H84").Value = 3278.0435
$ws.Range("J84").Value = 5519.6
$ws.Range("L84").Value = 55196
$ws.Range("N84").Value = -65804
